$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.961429476737976
$ws.Range("B1").Value = 2.085790157318115
$ws.Range("C1").Value = 2.086690187454224
$ws.Range("D1").Value = 2.598159313201904
$ws.Range("E1").Value = 3.367278099060059
